$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New C:H values for data rows 2-21 (shifted down by two rows from the
# previous data, with two brand-new rows of sensor readings inserted at
# the top of the "falling" sample window).
$data = @(
    @(-3.747647881507874, 4.277600646018982, 0.2108629420399666, 0.0218384321779012, -0.036499198526144, -0.0226020142436027),
    @(-3.872398495674133, 4.347799897193909, 0.4252039864659314, -0.00137444678694, 0.0125227374956011, 0.0018325957935303),
    @(-3.904903650283813, 4.009042739868163, 1.121196746826173, 0.064446285367012, 0.1020144969224929, 0.090408056974411),
    @(-4.147337973117831, 3.341572999954222, 1.590219795703887, -0.0655152946710586, 0.2060143053531646, 0.2438879460096359),
    @(-4.578064322471619, 2.440073847770686, 1.919119298458106, 0.1971567571163177, 0.3791945874691009, -0.0610865242779254),
    @(-4.763241469860086, 0.9537072181701609, 3.53445836901665, -0.3149010241031647, 0.4572326242923736, 0.2616030275821686),
    @(-6.981464385986306, -0.2511000633239622, 4.341212272644047, -0.2191479057073593, 1.047023057937622, -0.107512280344963),
    @(-1.964182496070859, 2.53690141439438, 5.374873697757723, 1.168585181236267, -2.228741884231567, -0.9682214260101318),
    @(0.6403729915618919, 8.453629255294809, -1.445600330829628, 0.5827654600143433, 0.8587238192558289, -2.957046031951904),
    @(25.52450680732733, 10.11635231971741, -29.67163440585143, -0.6496551632881165, -2.849991798400879, 1.479821085929871),
    @(21.21333217620842, 6.795180916786184, -24.62865638732902, 0.3932445049285888, -0.2875648140907287, -0.4578435122966766),
    @(1.512938618659968, 4.070763558149341, -2.205098286271084, 0.0209221355617046, 0.8381071090698242, -0.4121813178062439),
    @(0.09360861778259269, 5.09270715713501, -0.006848305463790894, -0.3355177342891693, -1.3600914478302, 0.2121229618787765),
    @(0.4487386941909907, 5.822665631771091, -1.949815332889539, 0.1876883506774902, -1.909106612205505, -0.6495024561882019),
    @(3.47012042999266, 6.646041870117178, 2.636732816696162, -0.6143777370452881, 0.0381790772080421, -0.1973094791173935),
    @(-0.02146565914155785, 4.178134083747851, 1.41586604714393, -0.2321287989616394, -0.2197587788105011, 0.07544185966253281),
    @(0.8478746414184657, 4.644729614257818, 0.6306794285774192, -0.051312681287527, -0.1487456858158111, 0.2446515262126922),
    @(0.5417299270629863, 5.275527954101563, -0.05354586243629561, 0.0394008085131645, 0.1130100712180137, 0.1357648074626922),
    @(-0.06840181350708066, 5.362720489501952, -0.7121562957763687, 0.038026362657547, 0.1304197311401367, 0.0140499006956815),
    @(-0.2478119134902955, 5.057902812957763, -0.9678008556365962, -0.0109955742955207, 0.0171042270958423, 0.008399397134780801)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 3).Value2 = $row[0]
    $ws.Cells.Item($r, 4).Value2 = $row[1]
    $ws.Cells.Item($r, 5).Value2 = $row[2]
    $ws.Cells.Item($r, 6).Value2 = $row[3]
    $ws.Cells.Item($r, 7).Value2 = $row[4]
    $ws.Cells.Item($r, 8).Value2 = $row[5]
}

# Remove the now-obsolete last data row (row 22), shrinking the sheet
# from A1:H22 down to A1:H21.
$ws.Rows.Item(22).Delete()
